{"js": "// Remove the trailing \"Requisitos\" section (its Heading2 title paragraph and\n// the following \"LOT2040 - Engenharia Gen\u00e9tica (Requisito fraco)\" list-bullet\n// paragraph) from the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text,style\");\nawait context.sync();\n\n// Locate the \"Requisitos\" heading paragraph.\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text.trim();\n  const style = (paragraphs.items[i].style || \"\").trim();\n  if (text === \"Requisitos\" && style === \"Heading 2\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex !== -1) {\n  // Delete every paragraph from \"Requisitos\" through the end of the body\n  // (i.e. the heading itself plus the requirement list that follows it).\n  const toDelete = [];\n  for (let i = startIndex; i < paragraphs.items.length; i++) {\n    toDelete.push(paragraphs.items[i]);\n  }\n  // Delete from last to first so indices/ranges stay valid as we go.\n  for (let i = toDelete.length - 1; i >= 0; i--) {\n    toDelete[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Requisitos\" section (its Heading 2 title paragraph and\n# the following \"LOT2040 - Engenharia Gen\u00e9tica (Requisito fraco)\" list-bullet\n# paragraph) from the end of the document.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Requisitos\" heading paragraph.\n$count = $d.Paragraphs.Count\n$startIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $txt = $p.Range.Text.Trim()\n    $styleName = $p.Style.NameLocal\n    if ($txt -eq \"Requisitos\" -and $styleName -eq \"Heading 2\") {\n        $startIndex = $i\n        break\n    }\n}\n\nif ($startIndex -ne -1) {\n    # Build a range spanning from the start of the \"Requisitos\" heading\n    # through the end of the document (the requirement list that follows it)\n    # and delete it in one shot.\n    $startPara = $d.Paragraphs.Item($startIndex)\n    $lastPara = $d.Paragraphs.Item($count)\n    $r = $d.Range($startPara.Range.Start, $lastPara.Range.End)\n    $r.Delete()\n}\n"}
